$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data to the Commodity table
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "wind_onshore"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "heat"

# Resize the table (ListObject) to include the new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:B6"))

# Update selection to mirror the saved state (next empty row)
$ws.Range("B7").Select()
